$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 54
$ws1.Range("F4").Value = 27
$ws1.Range("F6").Value = 38
$ws1.Range("F8").Value = 512
$ws1.Range("F9").Value = 3639
$ws1.Range("F10").Value = 61

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 4

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 54
$ws4.Range("F4").Value = 27
$ws4.Range("F6").Value = 38
$ws4.Range("F8").Value = 512
$ws4.Range("F9").Value = 3639
$ws4.Range("F10").Value = 61
$ws4.Range("F12").Value = 4
